$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/20/2024  Through  5/26/2024"

# --- Column H width widened to match column E (best-fit growth) ---
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# --- Cells that flip from a numeric value to the literal text "0" placeholder ---
# (copy an existing text-"0" cell, e.g. C14, to also bring along its style s="14")
$ws.Cells.Item(14, 3).Copy($ws.Cells.Item(15, 3))
$ws.Cells.Item(14, 3).Copy($ws.Cells.Item(27, 3))
$ws.Cells.Item(14, 3).Copy($ws.Cells.Item(29, 3))
$ws.Cells.Item(14, 3).Copy($ws.Cells.Item(30, 3))
$ws.Cells.Item(14, 3).Copy($ws.Cells.Item(31, 6))

# --- C28 flips from the text "0" placeholder back to a real number; pick up the
#     numeric style (s="15") from its neighbour D28 first, then set the value ---
$ws.Cells.Item(28, 4).Copy($ws.Cells.Item(28, 3))
$ws.Cells.Item(28, 3).Value = 3

# --- Remaining plain numeric value updates (style unchanged) ---
$ws.Cells.Item(15, 5).Value = -100
$ws.Cells.Item(15, 6).Value = 2
$ws.Cells.Item(15, 8).Value = -60
$ws.Cells.Item(15, 10).Value = 17
$ws.Cells.Item(15, 11).Value = -23.529411764705
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).Value = -27.777777777777
$ws.Cells.Item(16, 3).Value = 8
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = 14.285714285714
$ws.Cells.Item(16, 6).Value = 27
$ws.Cells.Item(16, 7).Value = 25
$ws.Cells.Item(16, 8).Value = 8
$ws.Cells.Item(16, 9).Value = 128
$ws.Cells.Item(16, 10).Value = 110
$ws.Cells.Item(16, 11).Value = 16.363636363636
$ws.Cells.Item(16, 12).Value = 36.170212765957
$ws.Cells.Item(16, 13).Value = 11.304347826087
$ws.Cells.Item(16, 14).Value = -82.772543741588
$ws.Cells.Item(17, 3).Value = 9
$ws.Cells.Item(17, 4).Value = 12
$ws.Cells.Item(17, 5).Value = -25
$ws.Cells.Item(17, 7).Value = 43
$ws.Cells.Item(17, 8).Value = -6.976744186046
$ws.Cells.Item(17, 9).Value = 223
$ws.Cells.Item(17, 10).Value = 190
$ws.Cells.Item(17, 11).Value = 17.368421052631
$ws.Cells.Item(17, 12).Value = 9.852216748768
$ws.Cells.Item(17, 13).Value = 120.792079207921
$ws.Cells.Item(17, 14).Value = -26.158940397351
$ws.Cells.Item(18, 3).Value = 6
$ws.Cells.Item(18, 4).Value = 7
$ws.Cells.Item(18, 5).Value = -14.285714285714
$ws.Cells.Item(18, 6).Value = 14
$ws.Cells.Item(18, 7).Value = 20
$ws.Cells.Item(18, 8).Value = -30
$ws.Cells.Item(18, 9).Value = 84
$ws.Cells.Item(18, 10).Value = 93
$ws.Cells.Item(18, 11).Value = -9.677419354838
$ws.Cells.Item(18, 12).Value = -19.230769230769
$ws.Cells.Item(18, 13).Value = -57.788944723618
$ws.Cells.Item(18, 14).Value = -89.5
$ws.Cells.Item(19, 3).Value = 9
$ws.Cells.Item(19, 4).Value = 19
$ws.Cells.Item(19, 5).Value = -52.631578947368
$ws.Cells.Item(19, 6).Value = 77
$ws.Cells.Item(19, 7).Value = 70
$ws.Cells.Item(19, 8).Value = 10
$ws.Cells.Item(19, 9).Value = 383
$ws.Cells.Item(19, 10).Value = 377
$ws.Cells.Item(19, 11).Value = 1.591511936339
$ws.Cells.Item(19, 12).Value = 13.988095238095
$ws.Cells.Item(19, 13).Value = 47.307692307692
$ws.Cells.Item(19, 14).Value = -13.54401805869
$ws.Cells.Item(20, 3).Value = 5
$ws.Cells.Item(20, 4).Value = 8
$ws.Cells.Item(20, 5).Value = -37.5
$ws.Cells.Item(20, 6).Value = 26
$ws.Cells.Item(20, 7).Value = 27
$ws.Cells.Item(20, 8).Value = -3.703703703703
$ws.Cells.Item(20, 9).Value = 139
$ws.Cells.Item(20, 10).Value = 130
$ws.Cells.Item(20, 11).Value = 6.923076923076
$ws.Cells.Item(20, 12).Value = 47.872340425531
$ws.Cells.Item(20, 13).Value = 31.132075471698
$ws.Cells.Item(20, 14).Value = -87.46618575293
$ws.Cells.Item(21, 3).Value = 37
$ws.Cells.Item(21, 4).Value = 54
$ws.Cells.Item(21, 5).Value = -31.481481481481
$ws.Cells.Item(21, 6).Value = 186
$ws.Cells.Item(21, 7).Value = 190
$ws.Cells.Item(21, 8).Value = -2.105263157894
$ws.Cells.Item(21, 9).Value = 970
$ws.Cells.Item(21, 10).Value = 920
$ws.Cells.Item(21, 11).Value = 5.434782608695
$ws.Cells.Item(21, 12).Value = 14.521841794569
$ws.Cells.Item(21, 13).Value = 22.320302648171
$ws.Cells.Item(21, 14).Value = -71.687098657326
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(22, 4).Value = 3
$ws.Cells.Item(22, 5).Value = -66.666666666666
$ws.Cells.Item(22, 6).Value = 3
$ws.Cells.Item(22, 7).Value = 4
$ws.Cells.Item(22, 8).Value = -25
$ws.Cells.Item(22, 9).Value = 12
$ws.Cells.Item(22, 10).Value = 15
$ws.Cells.Item(22, 11).Value = -20
$ws.Cells.Item(22, 12).Value = 100
$ws.Cells.Item(22, 13).Value = -7.692307692307
$ws.Cells.Item(23, 3).Value = 5
$ws.Cells.Item(23, 4).Value = 5
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 28
$ws.Cells.Item(23, 8).Value = 7.692307692307
$ws.Cells.Item(23, 9).Value = 87
$ws.Cells.Item(23, 10).Value = 96
$ws.Cells.Item(23, 11).Value = -9.375
$ws.Cells.Item(23, 12).Value = 2.35294117647
$ws.Cells.Item(23, 13).Value = 64.150943396226
$ws.Cells.Item(24, 3).Value = 52
$ws.Cells.Item(24, 4).Value = 67
$ws.Cells.Item(24, 5).Value = -22.388059701492
$ws.Cells.Item(24, 6).Value = 190
$ws.Cells.Item(24, 7).Value = 196
$ws.Cells.Item(24, 8).Value = -3.061224489795
$ws.Cells.Item(24, 9).Value = 911
$ws.Cells.Item(24, 10).Value = 950
$ws.Cells.Item(24, 11).Value = -4.105263157894
$ws.Cells.Item(24, 12).Value = -11.121951219512
$ws.Cells.Item(24, 13).Value = 59.824561403508
$ws.Cells.Item(25, 3).Value = 29
$ws.Cells.Item(25, 4).Value = 38
$ws.Cells.Item(25, 5).Value = -23.684210526315
$ws.Cells.Item(25, 6).Value = 96
$ws.Cells.Item(25, 7).Value = 96
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 482
$ws.Cells.Item(25, 10).Value = 470
$ws.Cells.Item(25, 11).Value = 2.553191489361
$ws.Cells.Item(25, 12).Value = 25.194805194805
$ws.Cells.Item(26, 3).Value = 13
$ws.Cells.Item(26, 4).Value = 17
$ws.Cells.Item(26, 5).Value = -23.529411764705
$ws.Cells.Item(26, 6).Value = 81
$ws.Cells.Item(26, 7).Value = 80
$ws.Cells.Item(26, 8).Value = 1.25
$ws.Cells.Item(26, 9).Value = 371
$ws.Cells.Item(26, 10).Value = 328
$ws.Cells.Item(26, 11).Value = 13.109756097561
$ws.Cells.Item(26, 12).Value = -1.851851851851
$ws.Cells.Item(26, 13).Value = 6.303724928366
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(27, 5).Value = -100
$ws.Cells.Item(27, 6).Value = 2
$ws.Cells.Item(27, 7).Value = 8
$ws.Cells.Item(27, 8).Value = -75
$ws.Cells.Item(27, 10).Value = 24
$ws.Cells.Item(27, 11).Value = -16.666666666666
$ws.Cells.Item(27, 12).Value = 5.263157894736
$ws.Cells.Item(28, 4).Value = 2
$ws.Cells.Item(28, 5).Value = 50
$ws.Cells.Item(28, 6).Value = 5
$ws.Cells.Item(28, 7).Value = 7
$ws.Cells.Item(28, 8).Value = -28.571428571428
$ws.Cells.Item(28, 9).Value = 31
$ws.Cells.Item(28, 10).Value = 38
$ws.Cells.Item(28, 11).Value = -18.421052631578
$ws.Cells.Item(28, 12).Value = 6.896551724137
$ws.Cells.Item(29, 5).Value = -100
$ws.Cells.Item(29, 7).Value = 3
$ws.Cells.Item(29, 8).Value = -66.666666666666
$ws.Cells.Item(29, 10).Value = 14
$ws.Cells.Item(29, 11).Value = -71.428571428571
$ws.Cells.Item(29, 13).Value = -60
$ws.Cells.Item(30, 5).Value = -100
$ws.Cells.Item(30, 7).Value = 3
$ws.Cells.Item(30, 8).Value = -66.666666666666
$ws.Cells.Item(30, 10).Value = 13
$ws.Cells.Item(30, 11).Value = -76.923076923076
$ws.Cells.Item(30, 13).Value = -62.5
$ws.Cells.Item(31, 4).Value = 4
$ws.Cells.Item(31, 7).Value = 9
$ws.Cells.Item(31, 8).Value = -100
$ws.Cells.Item(31, 10).Value = 13
$ws.Cells.Item(31, 11).Value = -69.230769230769
